# Apply the edits described by the commit:
#  - Add an invigilator-availability constraint effectively changes the
#    duration/end_time/fill-colour-driving hour values for two existing
#    exams (rows 72 & 73), and adds a brand new exam (Exam75) as row 76.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1day_exam_venues")

# --- Update existing rows 72 and 73: duration (D) 9 -> 13, which bumps
#     end_time (F, =SUM(D,E)) from 11 -> 15, and the fill-colour-driving
#     hour column (H) from 16 -> 18.
$ws.Range("D72").Value = 13
$ws.Range("H72").Value = 18

$ws.Range("D73").Value = 13
$ws.Range("H73").Value = 18

# --- Append the new exam (Exam75) as row 76 ---
$newRow = 76
$ws.Range("A" + $newRow).Formula = "=ROW()-2"
$ws.Range("B" + $newRow).Value = "Exam75"
$ws.Range("C" + $newRow).Value2 = $ws.Range("C75").Value2
$ws.Range("D" + $newRow).Value = 13
$ws.Range("E" + $newRow).Value = 2
$ws.Range("F" + $newRow).Formula = "=SUM(D76,E76)"
$ws.Range("G" + $newRow).Value = 46
$ws.Range("H" + $newRow).Value = 18

# Match formatting (number format / style) of the preceding row for the
# new row's date & hour cells
$ws.Range("C75").Copy()
$ws.Range("C76").PasteSpecial(-4122) | Out-Null
$ws.Range("H75").Copy()
$ws.Range("H76").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update dimension / view bookkeeping to mirror Excel's own behaviour ---
$ws.Range("G76").Select() | Out-Null
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 57
